# Updates the "cryptos" price/volume table with freshly scraped values
# (GitHub Actions refresh job). Columns: B=Coin, C=Link, D=Price, E=Volume(1h).
# A leading "'" forces a numeric-looking price to stay text (matches the
# original inlineStr cell type instead of Excel auto-converting it to a
# number and dropping formatting like trailing zeros / thousand dots).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.194.19"
$ws.Range("E2").Value = "  +6.02%  "
$ws.Range("D3").Value = "2.586.10"
$ws.Range("E3").Value = "  +5.66%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "'588.69"
$ws.Range("E5").Value = "  +3.55%  "
$ws.Range("D6").Value = "'155.75"
$ws.Range("E6").Value = "  +6.86%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.48%  "
$ws.Range("D8").Value = "'0.544"
$ws.Range("E8").Value = "  +2.48%  "
$ws.Range("D9").Value = "2.611.99"
$ws.Range("E9").Value = "  +6.78%  "
$ws.Range("D10").Value = "'0.116"
$ws.Range("E10").Value = "  +4.67%  "
$ws.Range("E11").Value = "  -1.57%  "
$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").Value = "'5.33"
$ws.Range("E12").Value = "  +1.77%  "
$ws.Range("B13").Value = "Cardano"
$ws.Range("C13").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D13").Value = "'0.362"
$ws.Range("E13").Value = "  +3.69%  "
$ws.Range("D14").Value = "'29.67"
$ws.Range("E14").Value = "  +3.72%  "
$ws.Range("D15").Value = "'0.0000185"
$ws.Range("E15").Value = "  +6.02%  "
$ws.Range("D16").Value = "3.062.32"
$ws.Range("E16").Value = "  +5.70%  "
$ws.Range("D17").Value = "65.783.24"
$ws.Range("E17").Value = "  +5.37%  "
$ws.Range("D18").Value = "2.617.23"
$ws.Range("E18").Value = "  +6.79%  "
$ws.Range("D19").Value = "'8.16"
$ws.Range("E19").Value = "  +3.91%  "
$ws.Range("D20").Value = "'11.20"
$ws.Range("E20").Value = "  +3.94%  "
$ws.Range("D21").Value = "'354.48"
$ws.Range("E21").Value = "  +9.88%  "
$ws.Range("D22").Value = "'4.35"
$ws.Range("E22").Value = "  +5.25%  "
$ws.Range("D23").Value = "'2.29"
$ws.Range("E23").Value = "  +6.18%  "
$ws.Range("D24").Value = "'0.997"
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").Value = "'10.25"
$ws.Range("E25").Value = "  +2.90%  "
$ws.Range("D26").Value = "'66.23"
$ws.Range("E26").Value = "  +1.55%  "
$ws.Range("D27").Value = "'642.96"
$ws.Range("E27").Value = "  +1.37%  "
$ws.Range("E28").Value = "  +11.67%  "
$ws.Range("D30").Value = "'1.51"
$ws.Range("E30").Value = "  +6.97%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").Value = "'8.27"
$ws.Range("E32").Value = "  +5.22%  "
$ws.Range("D33").Value = "'1.91"
$ws.Range("E33").Value = "  +5.45%  "
$ws.Range("D34").Value = "'0.141"
$ws.Range("E34").Value = "  +7.47%  "
$ws.Range("D35").Value = "'1.62"
$ws.Range("E35").Value = "  +8.34%  "
$ws.Range("D36").Value = "'0.994"
$ws.Range("E36").Value = "  -0.48%  "
$ws.Range("D37").Value = "'5.00"
$ws.Range("E37").Value = "  +6.15%  "
$ws.Range("D38").Value = "'5.70"
$ws.Range("E38").Value = "  +8.73%  "
$ws.Range("D39").Value = "'2.95"
$ws.Range("E39").Value = "  +10.56%  "
$ws.Range("D40").Value = "'19.37"
$ws.Range("E40").Value = "  +4.97%  "
$ws.Range("D41").Value = "'156.22"
$ws.Range("E41").Value = "  +3.38%  "
$ws.Range("D42").Value = "'0.375"
$ws.Range("E42").Value = "  +2.90%  "
$ws.Range("E43").Value = "  +7.84%  "
$ws.Range("D44").Value = "'42.17"
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("D45").Value = "'163.89"
$ws.Range("E45").Value = "  +7.85%  "
$ws.Range("D46").Value = "0.0₆0311"
$ws.Range("E46").Value = "  +1.70%  "
$ws.Range("D47").Value = "'0.998"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("D48").Value = "'16.06"
$ws.Range("E48").Value = "  +4.96%  "
$ws.Range("E49").Value = "  +6.40%  "
$ws.Range("D50").Value = "'21.90"
$ws.Range("E50").Value = "  +9.53%  "
$ws.Range("D51").Value = "'0.639"
$ws.Range("E51").Value = "  +6.13%  "
